$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.471.83'
$ws.Range("E2").Value = '  -1.40%  '

$ws.Range("D3").Value = '''1.847.22'
$ws.Range("E3").Value = '  -1.08%  '

$ws.Range("D4").Value = '''0.9996'
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '''265.09'
$ws.Range("E5").Value = '  -2.78%  '

$ws.Range("D6").Value = '''0.9987'
$ws.Range("E6").Value = '  -0.08%  '

$ws.Range("D7").Value = '''0.5201'
$ws.Range("E7").Value = '  -1.79%  '

$ws.Range("D8").Value = '''0.3282'
$ws.Range("E8").Value = '  -3.01%  '

$ws.Range("D9").Value = '''0.06807'
$ws.Range("E9").Value = '  -0.09%  '

$ws.Range("D10").Value = '''18.90'
$ws.Range("E10").Value = '  -4.97%  '

$ws.Range("D11").Value = '''0.7805'
$ws.Range("E11").Value = '  -1.63%  '

$ws.Range("D12").Value = '''0.07767'
$ws.Range("E12").Value = '  +0.26%  '

$ws.Range("D13").Value = '''1.850.61'
$ws.Range("E13").Value = '  -0.95%  '

$ws.Range("D14").Value = '''88.16'
$ws.Range("E14").Value = '  -2.43%  '

$ws.Range("D15").Value = '''5.021'
$ws.Range("E15").Value = '  -2.14%  '

$ws.Range("D16").Value = '''0.9994'
$ws.Range("E16").Value = '  +0.07%  '

$ws.Range("D17").Value = '''13.94'
$ws.Range("E17").Value = '  -3.20%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '''0.000007986'
$ws.Range("E18").Value = '  -0.22%  '

$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").Value = '''1.000'
$ws.Range("E19").Value = '  +0.06%  '

$ws.Range("D20").Value = '''26.497.75'
$ws.Range("E20").Value = '  -1.42%  '

$ws.Range("D21").Value = '''2.074.34'
$ws.Range("E21").Value = '  -1.93%  '

$ws.Range("D22").Value = '''4.637'
$ws.Range("E22").Value = '  -1.21%  '

$ws.Range("D23").Value = '''9.590'
$ws.Range("E23").Value = '  -3.50%  '

$ws.Range("D24").Value = '''6.005'
$ws.Range("E24").Value = '  -1.04%  '

$ws.Range("D25").Value = '''144.39'
$ws.Range("E25").Value = '  -0.96%  '

$ws.Range("D26").Value = '''2.196'
$ws.Range("E26").Value = '  -7.98%  '

$ws.Range("D27").Value = '''1.669'
$ws.Range("E27").Value = '  +0.70%  '

$ws.Range("D28").Value = '''17.04'
$ws.Range("E28").Value = '  -1.10%  '

$ws.Range("D29").Value = '''112.46'
$ws.Range("E29").Value = '  -0.14%  '

$ws.Range("D30").Value = '''4.178'
$ws.Range("E30").Value = '  -3.43%  '

$ws.Range("D31").Value = '''4.142'
$ws.Range("E31").Value = '  -3.77%  '

$ws.Range("D32").Value = '''0.08731'
$ws.Range("E32").Value = '  -1.37%  '

$ws.Range("D33").Value = '''0.04839'
$ws.Range("E33").Value = '  -2.16%  '

$ws.Range("D34").Value = '''0.7251'
$ws.Range("E34").Value = '  -0.36%  '

$ws.Range("D35").Value = '''1.135'
$ws.Range("E35").Value = '  -2.38%  '

$ws.Range("D36").Value = '''2.851'
$ws.Range("E36").Value = '  -0.95%  '

$ws.Range("D37").Value = '''3.098'
$ws.Range("E37").Value = '  -3.02%  '

$ws.Range("D38").Value = '''0.01783'
$ws.Range("E38").Value = '  -3.39%  '

$ws.Range("D39").Value = '''2.226'
$ws.Range("E39").Value = '  -3.56%  '

$ws.Range("D40").Value = '''0.4895'
$ws.Range("E40").Value = '  -3.81%  '

$ws.Range("D41").Value = '''0.9141'
$ws.Range("E41").Value = '  -2.26%  '

$ws.Range("D42").Value = '''111.49'
$ws.Range("E42").Value = '  -4.01%  '

$ws.Range("D43").Value = '''6.090'
$ws.Range("E43").Value = '  -0.81%  '

$ws.Range("D44").Value = '''0.9989'
$ws.Range("E44").Value = '  -0.02%  '

$ws.Range("D45").Value = '''7.762'
$ws.Range("E45").Value = '  -3.18%  '

$ws.Range("D46").Value = '''0.4193'
$ws.Range("E46").Value = '  -5.14%  '

$ws.Range("D47").Value = '''0.05947'
$ws.Range("E47").Value = '  +0.16%  '

$ws.Range("D48").Value = '''9.140'
$ws.Range("E48").Value = '  -1.98%  '

$ws.Range("D49").Value = '''0.1249'
$ws.Range("E49").Value = '  -5.88%  '

$ws.Range("D50").Value = '''35.04'

$ws.Range("D51").Value = '''0.8892'
$ws.Range("E51").Value = '  +1.25%  '
